$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.111.27"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "2.520.39"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.39%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("D9").Value = "2.518.68"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("E11").Value = "  -2.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.348"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.46%  "
$ws.Range("D14").Value = "2.966.79"
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.02"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.14%  "
$ws.Range("D16").Value = "59.026.07"
$ws.Range("E16").Value = "  -0.47%  "
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("D18").Value = "2.519.67"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("E20").Value = "  -0.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.92%  "
$ws.Range("E24").Value = "  +3.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.423"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.42%  "
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.57%  "
$ws.Range("E29").Value = "  -3.68%  "
$ws.Range("D30").Value = "0.0₃0769"
$ws.Range("E30").Value = "  -1.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("E32").Value = "  +5.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "163.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.21%  "
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("E37").Value = "  -3.74%  "
$ws.Range("E38").Value = "  -2.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.57"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "286.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.04%  "
$ws.Range("E43").Value = "  -1.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "132.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.606"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0931"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("E49").Value = "  -0.96%  "
$ws.Range("E50").Value = "  -1.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.76%  "
